$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.683.53"
$ws.Range("E2").Value = "  +4.32%  "
$ws.Range("D3").Value = "1.758.86"
$ws.Range("E3").Value = "  +2.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.56"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.89%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9963"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4837"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.57%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2647"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.20%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06188"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Value = "1.751.13"
$ws.Range("E10").Value = "  +1.71%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "16.25"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.71%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.06952"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.10%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6104"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.87%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.529"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.43%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.75"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.47%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9960"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.02%  "
$ws.Range("D17").Value = "27.671.70"
$ws.Range("E17").Value = "  +4.96%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9940"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.13%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007104"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.42%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.58"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.16%  "
$ws.Range("D21").Value = "1.971.01"
$ws.Range("E21").Value = "  +1.39%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.504"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.19%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.478"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.55%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.140"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.97%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "141.37"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.46%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.872"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +7.16%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.39"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.94%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "109.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.384"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.94%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.987"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.92%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08078"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.80%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.711"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.25%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04671"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.18%  "
$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.024"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.47%  "
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.610"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.41%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6253"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.12%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9293"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.65%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.565"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.28%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.046"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.87%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9946"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.11%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01508"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.18%  "
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.720"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.49%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.66"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3878"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.31%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.948"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.12%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1165"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.50%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05365"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.09%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.894"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.90%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "30.13"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.33%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.257"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.58%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "51.47"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.19%  "
